$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, pushing the existing data
# (old rows 2-10) down to rows 3-11, and populate the new row 2 with
# the latest weekly price record.
$ws.Rows("2:2").Insert()

# The inserted row inherits the header row's formatting; clear it so the
# new row matches the plain formatting used by the rest of the data rows.
$ws.Range("A2:R2").ClearFormats()
# Restore the date/time number format on column D, matching the other
# data rows in this column.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 44473
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 300000000
$ws.Range("G2").Value = "Espárragos"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 2000
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Provincia de Linares"
$ws.Range("P2").Value = 2000
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
